$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D6").Value = 7
$ws.Range("E6").Value = 2
$ws.Range("H6").Value = 3
$ws.Range("I6").Value = 2
$ws.Range("F6").Select()
